$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update "Riders" (C) and "Average" (D) columns with new Madigan bike hours
$ws.Range("C2").Value = 269
$ws.Range("D2").Value = 269

$ws.Range("C3").Value = 260
$ws.Range("D3").Value = 260

$ws.Range("C4").Value = 270
$ws.Range("D4").Value = 270

$ws.Range("C5").Value = 266
$ws.Range("D5").Value = 266

$ws.Range("C6").Value = 234
$ws.Range("D6").Value = 257.5

$ws.Range("C7").Value = 96

$ws.Range("C8").Value = 72
$ws.Range("D8").Value = 78
